# Update header row of the active sheet so that the plan-import templates
# no longer rely on "sticky id" / "level" naming, switching to the new
# row-id / task / start-date / end-date headers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("PV-Test-03-t04-start-date")

$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Make sure this sheet is active, then move the selected/active cell.
$ws.Activate()
$ws.Range("F2").Select()
